$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 13 new match rows (rows 3-15) for Nicholas Pooran,
# keeping every value as text (matches the original sheet, which
# stores numeric-looking values like runs/balls/SR as text too).
$newRows = @(
    ,@(" Oct 30 2020", " Abu Dhabi", "Royals won by 7 wickets (with 15 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Nicholas Pooran ", "22", "10", "0", "3", "220.00")
    ,@(" Sep 27 2020", " Sharjah", "Royals won by 4 wickets (with 3 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Nicholas Pooran ", "25", "8", "1", "3", "312.50")
    ,@(" Oct 15 2020", " Sharjah", "Kings XI won by 8 wickets", "Kings XI Punjab", "Royal Challengers Bangalore", "Nicholas Pooran ", "6", "1", "0", "1", "600.00")
    ,@(" Oct 20 2020", " Dubai (DSC)", "Kings XI won by 5 wickets (with 6 balls remaining)", "Kings XI Punjab", "Delhi Capitals", "Nicholas Pooran ", "53", "28", "6", "3", "189.28")
    ,@(" Sep 24 2020", " Dubai (DSC)", "Kings XI won by 97 runs", "Kings XI Punjab", "Royal Challengers Bangalore", "Nicholas Pooran ", "17", "18", "1", "0", "94.44")
    ,@(" Oct 18 2020", " Dubai (DSC)", "Match tied (Kings XI won the one-over eliminator)", "Kings XI Punjab", "Mumbai Indians", "Nicholas Pooran ", "24", "12", "2", "2", "200.00")
    ,@(" Oct 4 2020", " Dubai (DSC)", "Super Kings won by 10 wickets (with 14 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Nicholas Pooran ", "33", "17", "1", "3", "194.11")
    ,@(" Oct 24 2020", " Dubai (DSC)", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Nicholas Pooran ", "32", "28", "2", "0", "114.28")
    ,@(" Oct 26 2020", " Sharjah", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kings XI Punjab", "Kolkata Knight Riders", "Nicholas Pooran ", "2", "3", "0", "0", "66.66")
    ,@(" Oct 10 2020", " Abu Dhabi", "KKR won by 2 runs", "Kings XI Punjab", "Kolkata Knight Riders", "Nicholas Pooran ", "16", "10", "2", "1", "160.00")
    ,@(" Oct 8 2020", " Dubai (DSC)", "Sunrisers won by 69 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Nicholas Pooran ", "77", "37", "5", "7", "208.10")
    ,@(" Sep 20 2020", " Dubai (DSC)", "Match tied (Capitals won the one-over eliminator)", "Kings XI Punjab", "Delhi Capitals", "Nicholas Pooran ", "0", "3", "0", "0", "0.00")
    ,@(" Oct 1 2020", " Abu Dhabi", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", "Nicholas Pooran ", "44", "27", "3", "2", "162.96")
)

$startRow = 3
$numRows = $newRows.Count
$numCols = 11

$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = $newRows[$i][$j]
    }
}

$endRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
# Force text storage (no automatic number coercion), write the block,
# then drop back to the default style so no extra cell format sticks.
$rng.NumberFormat = "@"
$rng.Value2 = $data
$rng.Style = "Normal"
